# Update generated output numbers (column F) on the "展览" and "全部类型" sheets
# to match a newer scrape run (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7923
$ws1.Range("F5").Value = 95
$ws1.Range("F17").Value = 5817
$ws1.Range("F18").Value = 175
$ws1.Range("F19").Value = 253
$ws1.Range("F20").Value = 1673
$ws1.Range("F21").Value = 235
$ws1.Range("F22").Value = 373

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 7923
$ws4.Range("F5").Value = 95
$ws4.Range("F18").Value = 5817
$ws4.Range("F20").Value = 175
$ws4.Range("F21").Value = 253
$ws4.Range("F22").Value = 1673
$ws4.Range("F23").Value = 235
$ws4.Range("F24").Value = 373

$wb.Save()
